# Thrust_Against_Mass_Calculations.xlsx - update experiment notes, PID values,
# chart marker/position, and active-sheet/zoom bookkeeping.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 - just a view zoom change (75% -> 95%)
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Select() | Out-Null
$excel.ActiveWindow.Zoom = 95

# ---------------------------------------------------------------------------
# Exp_data - new trial notes in column D, updated note text, view changes
# ---------------------------------------------------------------------------
$expData = $wb.Worksheets.Item("Exp_data")
$expData.Select() | Out-Null
$excel.ActiveWindow.Zoom = 95

$expData.Range("D2").Value = "6 min 10 seconds (370 seconds). 2Nd trial: 6 min 20 secs. 55% thrust for hovering at ground level"
$expData.Range("D4").Value = "5 min 59 seconds (353 seconds). Second trial: 6 min 19 seconds (Same thrust needed). 55% thrust for hovering above ground level"
$expData.Range("D6").Value = "5 min 58 seconds. 55.25% thrust"
$expData.Range("D7").Value = "5 min 58 seconds. 55.75% thrust"
$expData.Range("D8").Value = "5 min 52 seconds (57.25% thrust)"
$expData.Range("D9").Value = "5 min 45 seconds (57.5% thrust)"
$expData.Range("D17").Value = "6 min 45 seconds (405 seconds). 2Nd trial 6 min 50 seconds. Thrust reqd for hovering above ground level: 50%"

$expData.Range("D10").Select() | Out-Null

# Chart: shrink marker size and move/resize the plot on the sheet
$chartObj = $expData.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.MarkerSize = 4

# Target position/size taken from the chart's EMU offset+extent (12700 EMU = 1 pt)
$chartObj.Left = 5646600 / 12700
$chartObj.Top = 3956040 / 12700
$chartObj.Width = 4325400 / 12700
$chartObj.Height = 3229560 / 12700

# ---------------------------------------------------------------------------
# PID_Values - updated K3/L3/N3 figures, becomes the active/selected sheet
# ---------------------------------------------------------------------------
$pidValues = $wb.Worksheets.Item("PID_Values")
$pidValues.Select() | Out-Null
$excel.ActiveWindow.Zoom = 95

$pidValues.Range("K3").Value = 2300
$pidValues.Range("L3").Value = 40
$pidValues.Range("N3").Value = 2000

$pidValues.Range("N4").Select() | Out-Null

# PID_Values ends up the active sheet (tabSelected) / activeTab in the workbook
$pidValues.Select() | Out-Null
